$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from Sheet1 to hoja1
$ws.Name = "hoja1"

# Update header row shared strings (B/C/D/E -> L/U/M/N)
$ws.Range("B1").Value = "L"
$ws.Range("C1").Value = "U"
$ws.Range("D1").Value = "M"
$ws.Range("E1").Value = "N"

# Replace the rounded sample values with the full-precision originals
$ws.Range("A2").Value = 0.4233164630979916
$ws.Range("B2").Value = 1.128516274850555
$ws.Range("C2").Value = 0.1431947987728295
$ws.Range("D2").Value = -1.41236453473537
$ws.Range("E2").Value = 0.9487793268361908
$ws.Range("A3").Value = 0.9889028572022224
$ws.Range("B3").Value = -1.822289682647148
$ws.Range("C3").Value = 0.8897230282691831
$ws.Range("D3").Value = -0.9967536973471972
$ws.Range("E3").Value = -0.951892827650184
$ws.Range("A4").Value = -0.06285832484010689
$ws.Range("B4").Value = -0.587966387259811
$ws.Range("C4").Value = -0.5436270185561874
$ws.Range("D4").Value = -0.3845884973810081
$ws.Range("E4").Value = 1.262242542462966
$ws.Range("A5").Value = 2.089130200928322
$ws.Range("B5").Value = -0.5002712275417545
$ws.Range("C5").Value = -1.090679569698155
$ws.Range("D5").Value = -0.3353026979806761
$ws.Range("E5").Value = -1.062143519259482
$ws.Range("A6").Value = -1.081332878046013
$ws.Range("B6").Value = -1.590270386148148
$ws.Range("C6").Value = 1.355749667425276
$ws.Range("D6").Value = -0.348654115996958
$ws.Range("E6").Value = -0.4652896269021505
$ws.Range("A7").Value = 0.0133058808526776
$ws.Range("B7").Value = 0.1767885348951624
$ws.Range("C7").Value = -1.237786393230702
$ws.Range("D7").Value = -0.4119096260920766
$ws.Range("E7").Value = 2.284914218842941
$ws.Range("A8").Value = -0.3063111050674059
$ws.Range("B8").Value = 2.023690233726497
$ws.Range("C8").Value = 0.91272230261041
$ws.Range("D8").Value = 0.4743388227124177
$ws.Range("E8").Value = 0.65670927708812
$ws.Range("A9").Value = 0.5625345749784322
$ws.Range("B9").Value = 0.0288412741344594
$ws.Range("C9").Value = -1.126961147145649
$ws.Range("D9").Value = 0.2464037075798208
$ws.Range("E9").Value = 1.175491455761141
$ws.Range("A10").Value = -0.6522870633031013
$ws.Range("B10").Value = 0.9822288853780358
$ws.Range("C10").Value = 1.177280972076537
$ws.Range("D10").Value = -0.7352919403387764
$ws.Range("E10").Value = 2.597402338175312
$ws.Range("A11").Value = 0.694994461699469
$ws.Range("B11").Value = -0.5415108223460028
$ws.Range("C11").Value = -0.06668840776494429
$ws.Range("D11").Value = 0.3350891121534909
$ws.Range("E11").Value = 0.3450399444500379
